$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "410,276,116,100,等12项"
$ws.Range("B6").Value = "100,276,152,724,等13项"
$ws.Range("C6").Value = "124,826,842,392,等10项"
$ws.Range("E6").Value = "124,528,710,458,等13项"
$ws.Range("B7").Value = "616,381,484,752,等103项"
$ws.Range("C7").Value = "490,516,458,608,等111项"
$ws.Range("D7").Value = "116,490,458,752,等107项"
$ws.Range("E7").Value = "360,826,704,381,等91项"
$ws.Range("F7").Value = "710,458,490,251,等103项"
